# Remove the "Marine Drive" Business Women roster entry (Sandy Heatherington /
# Rhiannon Charles / marinedrive.bizwomen@gmail.com), which occupied the two
# sub-rows 21:22 under the merged "Marine Drive" block (A19:A22) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the two rows; remaining rows shift up and the A19:A22 merge
# collapses down to A19:A20 automatically.
$ws.Range("21:22").Delete()

# Re-point the print area, which used to cover the sheet through row 43;
# after removing 2 rows the used range now ends at row 41.
$ws.PageSetup.PrintArea = "`$A`$1:`$E`$41"

# Restore the active selection to the Marine Drive block.
$null = $ws.Range("A19:A20").Select()
